# custom accuracy + 데이터 1000개
# Round the last data row (row 5) to 2 decimal places, then drop the
# extra sample row (row 6) that is no longer part of the trimmed preview.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(7.21, 5.07, 0.8, 15.63, 12.44, 5.67, 26.4, 8.73, 3.79, 5.44, 6.27, 6.56, 1.82, 5.64, 7.97, 4.94, 0.73, 0.45, 78.54, 15.96, 5.21, 10.54, 5.45, 1.07, 12.25, 4.6, 4.19, 4.91, 6.54, 0.55, 24.29, 2.83, 6.51)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $values[$i]
}

$ws.Rows(6).Delete()
